$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): course columns E:L reordered, "Violence Response" column dropped ---
$ws.Range("E1").Value = "DHA Accommodations (1 hr)"
$ws.Range("F1").Value = "Leadership Training (4 hrs)"
$ws.Range("G1").Value = "MHS Customer Service (1 hr)"
$ws.Range("H1").Value = "Counterintelligence (1 hr)"
$ws.Range("I1").Value = "HIPAA Training (1 hr)"
$ws.Range("J1").Value = "Supervisor Safety Training (2 hrs)"
$ws.Range("K1").Value = "Employee Safety (1 hr)"
$ws.Range("L1").Value = "RandomCourse"

# --- Row 2: now John Doe, identified by numeric ID 1234567891 (previously parsed by email) ---
$ws.Range("B2").Value = 1234567891
$ws.Range("C2").Value = "John"
$ws.Range("D2").Value = "Doe"
$ws.Range("E2").Value = "Completed"
$ws.Range("F2").Value = "Completed"
$ws.Range("G2").Value = "Completed"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""

# --- Row 3: now Andrew Hartmann, identified by numeric ID 3322111234 ---
$ws.Range("B3").Value = 3322111234
$ws.Range("C3").Value = "Andrew"
$ws.Range("D3").Value = "Hartmann"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "NOT Completed"
$ws.Range("I3").Value = "NOT Completed"
$ws.Range("J3").Value = "LATE (completed)"
$ws.Range("K3").Value = "LATE (completed)"
$ws.Range("L3").Value = ""

# --- Row 4: now Nicholas Fletcher, identified by numeric ID 4453245321 ---
$ws.Range("B4").Value = 4453245321
$ws.Range("C4").Value = "Nicholas"
$ws.Range("D4").Value = "Fletcher"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = "Completed"

# --- Row 5: now J Cena, identified by numeric ID 5555555555 ---
$ws.Range("B5").Value = 5555555555
$ws.Range("C5").Value = "J"
$ws.Range("D5").Value = "Cena"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = "LATE (completed)"

# --- Drop column M entirely: dimension shrinks from A1:M5 to A1:L5 ---
$ws.Range("M1:M5").Delete()
